$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (empty paragraph right before
#    the "Search UI" heading). Word automatically renumbers the
#    remaining bookmark ids in document order when one is removed.
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ------------------------------------------------------------------
# 2. Update the "Primary Data" section paragraph text and re-insert
#    the "_GoBack" bookmark at its new location inside that text.
# ------------------------------------------------------------------
$oldTail = " You can download selected data in an Excel file, comma separated in a CSV file, or tab separated in a TXT file."
$newTail = " A User can download selected data in an Excel file, comma separated in a CSV file, or tab separated in a TXT or TSV file. Also a user can filter and sort the data before and download only a subset of the dataset. For filtering  use the funnel button next to the variable."

$found = $d.Content.Find.Execute($oldTail, $true, $true, $false, $false, $false, $true, 1, $false, $newTail, 2)

# Re-insert the _GoBack bookmark right after "filtering " (before the
# double space and "use the funnel button...").
$marker = "filtering  use the funnel button next to the variable."
$rng = $d.Content
$rng.Find.Execute($marker, $true, $true) | Out-Null
$bmPoint = $d.Range($rng.Start + 10, $rng.Start + 10)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# ------------------------------------------------------------------
# 3. The footer's PAGE field cached result changes from "6" to "5".
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$footer.Range.Find.Execute("6", $false, $false, $false, $false, $false, $true, 1, $false, "5", 2) | Out-Null
